$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted right after row 120 (old data shifts down
# by one row, old row 129 becomes new row 130). The new row 121 duplicates
# the (old) row 120 values, and row 120 itself gets updated with the new
# reporting date.

# 1. Insert a new row before the old row 121, shifting rows 121-129 down to 122-130.
$ws.Rows.Item(121).Insert()

# 2. Copy the full contents of row 120 (the most recent existing record)
#    into the newly inserted row 121.
$ws.Range("A120:T120").Copy()
$ws.Range("A121:T121").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0

# 3. Update row 120 with the new reporting date for this week's entry.
$ws.Range("D120").Value = 45223
